$d = $word.ActiveDocument

function Merge-Run($doc, $text) {
    # Locate the (unique) occurrence of $text and re-assert the very same
    # text onto the found Range. Word consolidates a Range spanning
    # several runs into a single run whenever the Range's Text is
    # reassigned to something that is a genuine change -- so we briefly
    # append a harmless sentinel character and then delete it again via
    # an explicit character-offset Range. (A plain Find.Execute(...,
    # Replace:=wdReplaceAll) would also merge the runs, but it pushes the
    # replacement text through the "smart quotes" AutoCorrect pass, which
    # would corrupt the straight apostrophe in one of the paragraphs.)
    $rng = $doc.Content
    [void]$rng.Find.Execute($text)
    $startPos = $rng.Start
    $rng.Text = $text + "#"
    $sentinel = $doc.Range($startPos + $text.Length, $startPos + $text.Length + 1)
    [void]$sentinel.Delete()
}

# 1) "Comandi" / " da linea di comando" / ":"  ->  single run
#    "Comandi da linea di comando:"
Merge-Run $d "Comandi da linea di comando:"

# 2) "Lancia il web server ... browser" / ")."  -> single run
Merge-Run $d "Lancia il web server di default di Angular (da chiamare nella cartella di un progetto Angular. Una volta partita potete controllare su localhost:4200. Se si salva il codice si vede l'app aggiornata sul browser)."

# 3) "Imposta il watch (c" / "ompila automaticamente il file " -> single run
Merge-Run $d "Imposta il watch (compila automaticamente il file "

# 4) " quando viene modificato" / ", " / "per uscire ctrl-C)" -> single run
Merge-Run $d " quando viene modificato, per uscire ctrl-C)"

# 5) " utile quando si devono " / "m" / "ettere molte configurazioni" -> single run,
#    then a brand-new run ". (Autocompila)" is appended right after it.
Merge-Run $d "utile quando si devono mettere molte configurazioni"

$rng = $d.Content
[void]$rng.Find.Execute("utile quando si devono mettere molte configurazioni")
$para = $rng.Paragraphs(1)
$prng = $para.Range
$prng.Collapse(0)
[void]$prng.MoveEnd(1, -1)
[void]$prng.InsertAfter(". (Autocompila)")

# 6) " ha visibilità solo nel blocco in cui si trova" / ", a differenza di " -> single run
Merge-Run $d " ha visibilità solo nel blocco in cui si trova, a differenza di "

Write-Output "done"
